# Insert a new data row at row 218 (shifting existing rows 218:262 down to
# 219:263) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("218").Insert()

$ws.Range("A218").Value = 9
$ws.Range("B218").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C218").Value = "Metropolitana"
$ws.Range("D218").Value = 44637
$ws.Range("E218").Value = 13
$ws.Range("F218").Value = 100112030
$ws.Range("G218").Value = "Poroto granado"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 70
$ws.Range("K218").Value = 24000
$ws.Range("L218").Value = 26000
$ws.Range("M218").Value = 24857
$ws.Range("N218").Value = "$/saco 25 kilos"
$ws.Range("O218").Value = "Provincia de Cardenal Caro"
$ws.Range("P218").Value = 994
$ws.Range("Q218").Value = 25
$ws.Range("R218").Value = "Hortaliza"
